$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete all rows whose tissueType (column A) equals "Breast", working from the
# bottom up so row indices of rows still to be processed are unaffected.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 2; $r--) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($val -eq "Breast") {
        $ws.Rows.Item($r).Delete()
    }
}

# Rename every remaining tissueType value ("Mammary gland" / "Mammary Gland")
# to "Mammary".
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -eq "Mammary gland" -or $val -eq "Mammary Gland") {
        $cell.Value = "Mammary"
    }
}

$ws.Range("A34").Select()
